$wb = $excel.ActiveWorkbook

# The "Dynamic" sheet is the active/first sheet containing the Tx_th header
$ws = $wb.Worksheets.Item("Dynamic")

# Rename header cell A1 from "Tx_th (deg)" to "Tx_el (deg)"
$ws.Range("A1").Value = "Tx_el (deg)"

# Update the selected cell shown in the sheet view from B16 to A9
$ws.Activate()
$ws.Range("A9").Select()
